# Update the "productAriaLabel" text for a handful of rows whose products
# went out of online stock: insert " - Online kein Bestand" before the
# trailing price text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M9").Value = "Ölz Vollkorn Sandwich Toast Soft - Online kein Bestand 4.35 Schweizer Franken"
$ws.Range("M230").Value = "Betty Bossi Kuchenteig -35% Fett - Online kein Bestand 2.10 Schweizer Franken"
$ws.Range("M384").Value = "Leisi Kuchenteig rund ausgewallt Ø32cm glutenfrei - Online kein Bestand 4.95 Schweizer Franken"

# Refresh the crawl timestamp (column O) for every data row (2-397) to
# reflect the later re-crawl time of this snapshot.
$newTimestamp = "2023-01-02 20:49:56"
for ($row = 2; $row -le 397; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
